$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B44 from text "2" to a real number 2
$ws.Range("B44").Value = 2

# Add new row 45 with data
$ws.Range("A45").Value = "Ying Tang"
$ws.Range("B45").Value = "'3"
$ws.Range("C45").Value = "无"
$ws.Range("D45").Value = "QSN"
$ws.Range("E45").Value = "THE"
$ws.Range("F45").Value = "77ff87fb-cfc5-44ac-a4b7-cb33b05fed6f"
$ws.Range("G45").Value = "ByQpn1ZA-_annotated.xlsx"
$ws.Range("H45").Value = "If we know the regularization is fundamentally and mathematically wrong, why do we investigate its performance?"
